$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell reference and its new value. All columns in
# this sheet are stored as text (even numeric-looking IDs/ranks/counts), so
# we force the '@' (text) number format before writing, then restore General.
$edits = @(
  @{ Cell = "A2"; Value = '57991' },
  @{ Cell = "E2"; Value = '2597' },
  @{ Cell = "A3"; Value = '62918' },
  @{ Cell = "E3"; Value = '2524' },
  @{ Cell = "A5"; Value = '44598' },
  @{ Cell = "A7"; Value = '60367' },
  @{ Cell = "E7"; Value = '2557' },
  @{ Cell = "A8"; Value = '8021' },
  @{ Cell = "B8"; Value = '53060417' },
  @{ Cell = "C8"; Value = '㊥老纳信耶稣' },
  @{ Cell = "E8"; Value = '4664' },
  @{ Cell = "A9"; Value = '11070' },
  @{ Cell = "B9"; Value = '49710892' },
  @{ Cell = "C9"; Value = 'MMMMMMM' },
  @{ Cell = "E9"; Value = '4514' },
  @{ Cell = "A10"; Value = '11694' },
  @{ Cell = "E10"; Value = '4487' },
  @{ Cell = "A11"; Value = '19154' },
  @{ Cell = "E11"; Value = '4207' },
  @{ Cell = "A12"; Value = '24369' },
  @{ Cell = "E12"; Value = '4052' },
  @{ Cell = "A13"; Value = '40373' },
  @{ Cell = "E13"; Value = '3351' },
  @{ Cell = "A14"; Value = '58425' },
  @{ Cell = "A15"; Value = '65211' },
  @{ Cell = "A17"; Value = '12571' },
  @{ Cell = "E17"; Value = '4446' },
  @{ Cell = "A18"; Value = '13709' },
  @{ Cell = "B18"; Value = '54698813' },
  @{ Cell = "C18"; Value = '閃亮唐老鴨' },
  @{ Cell = "E18"; Value = '4397' },
  @{ Cell = "A19"; Value = '14210' },
  @{ Cell = "B19"; Value = '31495601' },
  @{ Cell = "C19"; Value = '陈晓军' },
  @{ Cell = "E19"; Value = '4378' },
  @{ Cell = "A20"; Value = '15896' },
  @{ Cell = "E20"; Value = '4312' },
  @{ Cell = "A21"; Value = '17452' },
  @{ Cell = "E21"; Value = '4261' },
  @{ Cell = "A22"; Value = '21091' },
  @{ Cell = "B22"; Value = '54085771' },
  @{ Cell = "C22"; Value = '㊥Matthieu' },
  @{ Cell = "E22"; Value = '4151' },
  @{ Cell = "A23"; Value = '21811' },
  @{ Cell = "B23"; Value = '55769051' },
  @{ Cell = "C23"; Value = '㊥叮叮当.' },
  @{ Cell = "E23"; Value = '4128' },
  @{ Cell = "A24"; Value = '23728' },
  @{ Cell = "B24"; Value = '56732705' },
  @{ Cell = "C24"; Value = '时间温柔皆遗憾' },
  @{ Cell = "E24"; Value = '4072' },
  @{ Cell = "A25"; Value = '27640' },
  @{ Cell = "E25"; Value = '3990' },
  @{ Cell = "A26"; Value = '28899' },
  @{ Cell = "E26"; Value = '3974' },
  @{ Cell = "A27"; Value = '30111' },
  @{ Cell = "E27"; Value = '3927' },
  @{ Cell = "A28"; Value = '31312' },
  @{ Cell = "E28"; Value = '3872' },
  @{ Cell = "A29"; Value = '40915' },
  @{ Cell = "E29"; Value = '3314' },
  @{ Cell = "A30"; Value = '42734' },
  @{ Cell = "E30"; Value = '3188' },
  @{ Cell = "A31"; Value = '556' },
  @{ Cell = "E31"; Value = '5347' },
  @{ Cell = "A32"; Value = '8481' },
  @{ Cell = "B32"; Value = '7852598' },
  @{ Cell = "C32"; Value = 'seiji' },
  @{ Cell = "E32"; Value = '4640' },
  @{ Cell = "A33"; Value = '9222' },
  @{ Cell = "E33"; Value = '4602' },
  @{ Cell = "A34"; Value = '10420' },
  @{ Cell = "B34"; Value = '11582001' },
  @{ Cell = "C34"; Value = 'iMinatoX4' },
  @{ Cell = "E34"; Value = '4545' },
  @{ Cell = "A35"; Value = '11430' },
  @{ Cell = "B35"; Value = '45967307' },
  @{ Cell = "C35"; Value = 'Ricky' },
  @{ Cell = "E35"; Value = '4498' },
  @{ Cell = "A36"; Value = '13265' },
  @{ Cell = "B36"; Value = '38995116' },
  @{ Cell = "C36"; Value = '"Ramesh Pavai Nam"' },
  @{ Cell = "E36"; Value = '4417' },
  @{ Cell = "A37"; Value = '13561' },
  @{ Cell = "B37"; Value = '56133764' },
  @{ Cell = "C37"; Value = 'ustcarter' },
  @{ Cell = "E37"; Value = '4404' },
  @{ Cell = "A38"; Value = '13644' },
  @{ Cell = "B38"; Value = '55317038' },
  @{ Cell = "C38"; Value = 'necman12345' },
  @{ Cell = "E38"; Value = '4400' },
  @{ Cell = "A39"; Value = '16449' },
  @{ Cell = "B39"; Value = '6809364' },
  @{ Cell = "C39"; Value = '"Scorp IP"' },
  @{ Cell = "E39"; Value = '4292' },
  @{ Cell = "A40"; Value = '16648' },
  @{ Cell = "E40"; Value = '4287' },
  @{ Cell = "A41"; Value = '17075' },
  @{ Cell = "B41"; Value = '26280580' },
  @{ Cell = "C41"; Value = '꧁SSS.TIGRESS꧂ᶻᵍˣ' },
  @{ Cell = "E41"; Value = '4274' },
  @{ Cell = "A42"; Value = '22471' },
  @{ Cell = "B42"; Value = '47459684' },
  @{ Cell = "C42"; Value = '㊥阿闹切克闹' },
  @{ Cell = "E42"; Value = '4108' },
  @{ Cell = "A43"; Value = '25384' },
  @{ Cell = "B43"; Value = '56379103' },
  @{ Cell = "C43"; Value = 'Globalking' },
  @{ Cell = "E43"; Value = '4022' },
  @{ Cell = "A44"; Value = '29487' },
  @{ Cell = "B44"; Value = '56573048' },
  @{ Cell = "C44"; Value = 'Xiaotian' },
  @{ Cell = "E44"; Value = '3953' },
  @{ Cell = "A45"; Value = '33696' },
  @{ Cell = "B45"; Value = '50837459' },
  @{ Cell = "C45"; Value = 'NINE日' },
  @{ Cell = "E45"; Value = '3754' },
  @{ Cell = "A46"; Value = '36066' },
  @{ Cell = "B46"; Value = '52997727' },
  @{ Cell = "C46"; Value = 'larios' },
  @{ Cell = "E46"; Value = '3625' },
  @{ Cell = "A47"; Value = '36144' },
  @{ Cell = "B47"; Value = '58203298' },
  @{ Cell = "C47"; Value = '权旨qua' },
  @{ Cell = "E47"; Value = '3621' },
  @{ Cell = "A48"; Value = '38896' },
  @{ Cell = "E48"; Value = '3450' },
  @{ Cell = "A49"; Value = '41231' },
  @{ Cell = "E49"; Value = '3294' },
  @{ Cell = "A50"; Value = '52258' },
  @{ Cell = "E50"; Value = '2740' },
  @{ Cell = "A51"; Value = '55077' },
  @{ Cell = "E51"; Value = '2660' },
  @{ Cell = "A52"; Value = '57366' },
  @{ Cell = "A53"; Value = '59304' },
  @{ Cell = "E53"; Value = '2574' },
  @{ Cell = "A54"; Value = '65320' },
  @{ Cell = "A55"; Value = '57841' },
  @{ Cell = "E55"; Value = '2600' },
  @{ Cell = "A56"; Value = '57081' },
  @{ Cell = "E56"; Value = '2614' },
  @{ Cell = "A59"; Value = '30570' },
  @{ Cell = "E59"; Value = '3907' },
  @{ Cell = "A60"; Value = '43704' },
  @{ Cell = "E60"; Value = '3121' },
  @{ Cell = "A61"; Value = '46834' },
  @{ Cell = "B61"; Value = '11645391' },
  @{ Cell = "C61"; Value = '"omar omar"' },
  @{ Cell = "E61"; Value = '2949' },
  @{ Cell = "A62"; Value = '57164' },
  @{ Cell = "B62"; Value = '55499394' },
  @{ Cell = "C62"; Value = 'Player-55499394' },
  @{ Cell = "E62"; Value = '2612' },
  @{ Cell = "A63"; Value = '66540' },
  @{ Cell = "E63"; Value = '2492' },
  @{ Cell = "A64"; Value = '93848' },
  @{ Cell = "E64"; Value = '1510' },
  @{ Cell = "A65"; Value = '108402' },
  @{ Cell = "E65"; Value = '1280' },
  @{ Cell = "A79"; Value = '50501' },
  @{ Cell = "E79"; Value = '2800' },
  @{ Cell = "A82"; Value = '115244' },
  @{ Cell = "E82"; Value = '1154' },
  @{ Cell = "A83"; Value = '128486' },
  @{ Cell = "B83"; Value = '58174442' },
  @{ Cell = "C83"; Value = 'Player-58174442' },
  @{ Cell = "E83"; Value = '1020' },
  @{ Cell = "A84"; Value = '159313' },
  @{ Cell = "B84"; Value = '15695258' },
  @{ Cell = "C84"; Value = 'Player-15695258' },
  @{ Cell = "E84"; Value = '1000' },
  @{ Cell = "B85"; Value = '58572199' },
  @{ Cell = "C85"; Value = '你干嘛～哎呦～' },
  @{ Cell = "B86"; Value = '58766144' },
  @{ Cell = "C86"; Value = 'EquablePrecedence38' },
  @{ Cell = "B87"; Value = '29355299' },
  @{ Cell = "C87"; Value = 'Player-29355299' },
  @{ Cell = "B88"; Value = '58910668' },
  @{ Cell = "C88"; Value = 'BrittleAuthor33' },
  @{ Cell = "B89"; Value = '55745105' },
  @{ Cell = "C89"; Value = 'eldeniz' }
)

foreach ($edit in $edits) {
  $target = $ws.Range($edit.Cell)
  $target.NumberFormat = "@"
  $target.Value = $edit.Value
  $target.NumberFormat = "General"
}
